$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.075.93"
$ws.Range("E2").Value = "  +1.97%  "

# Row 3
$ws.Range("D3").Value = "3.808.98"
$ws.Range("E3").Value = "  +0.57%  "

# Row 4
$ws.Range("D4").Value = "'0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "

# Row 5
$ws.Range("D5").Value = "'630.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.20%  "

# Row 6
$ws.Range("D6").Value = "'164.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "

# Row 7
$ws.Range("D7").Value = "3.806.83"
$ws.Range("E7").Value = "  +0.56%  "

# Row 8
$ws.Range("E8").Value = "  +0.32%  "

# Row 9
$ws.Range("D9").Value = "'0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.74%  "

# Row 10
$ws.Range("E10").Value = "  +2.05%  "

# Row 11
$ws.Range("D11").Value = "'0.453"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.65%  "

# Row 12
$ws.Range("D12").Value = "'6.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.29%  "

# Row 13
$ws.Range("D13").Value = "'0.0000249"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.20%  "

# Row 14
$ws.Range("D14").Value = "'35.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.57%  "

# Row 15
$ws.Range("D15").Value = "4.448.23"
$ws.Range("E15").Value = "  +0.63%  "

# Row 16
$ws.Range("D16").Value = "3.860.99"
$ws.Range("E16").Value = "  +2.32%  "

# Row 17
$ws.Range("D17").Value = "69.048.46"
$ws.Range("E17").Value = "  +1.97%  "

# Row 18
$ws.Range("D18").Value = "'17.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.08%  "

# Row 19
$ws.Range("D19").Value = "'7.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.81%  "

# Row 20
$ws.Range("E20").Value = "  -0.16%  "

# Row 21
$ws.Range("D21").Value = "'465.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.01%  "

# Row 22
$ws.Range("D22").Value = "'9.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.96%  "

# Row 23
$ws.Range("D23").Value = "'0.707"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.73%  "

# Row 24
$ws.Range("D24").Value = "'0.0000150"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.14%  "

# Row 25
$ws.Range("D25").Value = "'83.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.08%  "

# Row 26
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'11.97"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "

# Row 27
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").Value = "'2.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.46%  "

# Row 28
$ws.Range("D28").Value = "'10.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.09%  "

# Row 29
$ws.Range("E29").Value = "  -0.02%  "

# Row 30
$ws.Range("D30").Value = "3.959.30"
$ws.Range("E30").Value = "  +0.61%  "

# Row 31
$ws.Range("E31").Value = "  +2.92%  "

# Row 32
$ws.Range("D32").Value = "'2.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.93%  "

# Row 33
$ws.Range("D33").Value = "'7.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.56%  "

# Row 34
$ws.Range("D34").Value = "'29.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "

# Row 35
$ws.Range("E35").Value = "  +0.09%  "

# Row 36
$ws.Range("D36").Value = "'9.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.72%  "

# Row 37
$ws.Range("E37").Value = "  +3.57%  "

# Row 38
$ws.Range("E38").Value = "  +7.50%  "

# Row 39
$ws.Range("D39").Value = "'3.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.60%  "

# Row 40
$ws.Range("D40").Value = "'5.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.09%  "

# Row 41
$ws.Range("D41").Value = "'0.974"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.30%  "

# Row 42
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.04%  "

# Row 43
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("D44").Value = "'156.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.19%  "

# Row 45
$ws.Range("D45").Value = "'0.299"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.70%  "

# Row 46
$ws.Range("E46").Value = "  +6.57%  "

# Row 47
$ws.Range("D47").Value = "'43.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.59%  "

# Row 48
$ws.Range("D48").Value = "'46.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.14%  "

# Row 49
$ws.Range("E49").Value = "  +2.99%  "

# Row 50
$ws.Range("D50").Value = "'8.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.24%  "

# Row 51
$ws.Range("D51").Value = "'0.000275"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.76%  "
